# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) previously held a different stat ("Strike#"); this
# re-derives/re-writes it with the new K values for each logged game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = [ordered]@{
    2  = 0
    3  = 2
    4  = 1
    5  = 0
    6  = 0
    7  = 2
    8  = 0
    9  = 2
    10 = 1
    11 = 3
    12 = 1
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 2
    18 = 2
    19 = 3
    20 = 0
    21 = 1
    22 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
